$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-17 12:50:50"
$wsZhCn.Range("H2").Value = "2016-03-17 12:51:10"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-17 12:50:54"
$wsDeDe.Range("H2").Value = "2016-03-17 12:51:16"
